# Atualização de bases das ligas, do dia: 18-02-2024 às 22:54
# This edit reorders a handful of match rows within the "Ecuador LigaPro Serie A"
# sheet: the row index in column A stays put, but the actual match data
# (columns B..AC: id, teams, odds, results, etc.) is permuted among a few rows,
# as if several rows had been re-sorted into a different (correct) order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each destination row, note which row currently holds the data that should
# end up there. Column A (the running counter) is left untouched in every row.
$rowMap = @{
    226 = 227
    227 = 226
    230 = 232
    231 = 233
    232 = 231
    233 = 230
    235 = 236
    236 = 235
    238 = 239
    239 = 241
    240 = 238
    241 = 240
}

# Snapshot the current content of columns B:AC for every row involved, before
# writing anything back - several rows are part of longer permutation cycles
# (e.g. 230 -> 232 -> 231 -> 233 -> 230), so the source data must be captured
# up-front rather than read-and-written one row at a time.
$snapshot = @{}
foreach ($row in $rowMap.Keys) {
    $snapshot[$row] = $ws.Range("B${row}:AC${row}").Value2
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $ws.Range("B${destRow}:AC${destRow}").Value2 = $snapshot[$srcRow]
}
